$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (old F/G shift right to G/H)
$ws.Columns.Item(6).Insert()

# New header cell F1 gets the "title" label
$ws.Cells.Item(1, 6).Value = "title"

# Match the new column's style/width to its left neighbor (column E)
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# Update the active selection to the new cell F8
$ws.Range("F8").Select()
